# ---------------------------------------------------------------------------
# 20220711 Readme NPOC/TN workbook update
#   - Adds a new "Dilution sheet" worksheet (after Sheet1) used to track
#     sample weights / DI dilution volumes before running NPOC/TN.
#   - Tweaks the "Total vol:" header on Sheet1 (trailing space + styling to
#     match the new sheet's header look).
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Build two tiny "prototype" cells off in an unused column of Sheet1 that we
# use purely to mint the two font looks we need (plain black Calibri, and
# bold black Calibri) plus the yellow header fill. We then propagate those
# looks onto the real cells with Copy/PasteSpecial(formats) so every cell
# that should look the same shares one style entry, and finally clear the
# scratch cells back out so they don't show up in anyone's used range.
# ---------------------------------------------------------------------------
$protoPlain = $ws1.Range("Z1")
$protoPlain.Value = "proto"
$protoPlain.Font.Name = "Calibri"
$protoPlain.Font.Color = 0

$protoBold = $ws1.Range("Z2")
$protoBold.Value = "proto"
$protoPlain.Copy()
$protoBold.PasteSpecial(-4122)
$protoBold.Font.Bold = $true

$protoBoldFill = $ws1.Range("Z3")
$protoBoldFill.Value = "proto"
$protoBold.Copy()
$protoBoldFill.PasteSpecial(-4122)
$protoBoldFill.Interior.PatternColor = 0
$protoBoldFill.Interior.Color = 10086143

# ---------------------------------------------------------------------------
# Sheet1: restyle the "Total vol:" header (F1) and its neighbour (E1) to the
# plain black-Calibri look, and give "Total vol:" a trailing space.
# ---------------------------------------------------------------------------
$ws1.Range("F1").Value = "Total vol: "

$protoPlain.Copy()
$ws1.Range("E1:F1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Scratch cells no longer needed.
# ---------------------------------------------------------------------------
$protoPlain.Clear()
$protoBold.Clear()
$protoBoldFill.Clear()

# ---------------------------------------------------------------------------
# New "Dilution sheet" worksheet, placed right after Sheet1.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Dilution sheet"

$ws2.Range("A1").Value = "Sample"
$ws2.Range("B1").Value = "Vial wt (g)"
$ws2.Range("C1").Value = "Vial wt+Sample (g)"
$ws2.Range("D1").Value = "Sample Wt (g)"
$ws2.Range("E1").Value = "DI added (mL)"
$ws2.Range("F1").Value = "Total vol (mL)"
$ws2.Range("G1").Value = "Vial wt after addition (g)"
$ws2.Range("J1").Value = "*two different glass thickness-> explains differing vial wts"

# Bold + yellow-fill header cells.
$protoBoldFill2 = $ws1.Range("Z3")
$protoBoldFill2.Value = "proto"
$protoBoldFill2.Font.Name = "Calibri"
$protoBoldFill2.Font.Color = 0
$protoBoldFill2.Font.Bold = $true
$protoBoldFill2.Interior.PatternColor = 0
$protoBoldFill2.Interior.Color = 10086143
$protoBoldFill2.Copy()
$ws2.Range("A1").PasteSpecial(-4122)
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("F1").PasteSpecial(-4122)
$protoBoldFill2.Clear()

# Bold, no-fill header cells.
$protoBold2 = $ws1.Range("Z1")
$protoBold2.Value = "proto"
$protoBold2.Font.Name = "Calibri"
$protoBold2.Font.Color = 0
$protoBold2.Font.Bold = $true
$protoBold2.Copy()
$ws2.Range("B1").PasteSpecial(-4122)
$ws2.Range("C1").PasteSpecial(-4122)
$ws2.Range("E1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)
$ws2.Range("H1").PasteSpecial(-4122)
$protoBold2.Clear()

# Plain black header look used for the footnote + its empty neighbour.
$protoPlain2 = $ws1.Range("Z2")
$protoPlain2.Value = "proto"
$protoPlain2.Font.Name = "Calibri"
$protoPlain2.Font.Color = 0
$protoPlain2.Copy()
$ws2.Range("I1").PasteSpecial(-4122)
$ws2.Range("J1").PasteSpecial(-4122)
$protoPlain2.Clear()

$ws2.Range("A1:J1").Select()

$ws1.Select()
$ws1.Range("B20").Select()
